# Adds a new "05dec2025" data column (inserted before the existing
# "26nov2025" column, i.e. at column F) to both the "crosstab" and "annot"
# sheets, shifting the later nov-2025 columns one place to the right.
# Also updates the "04dec2025" (E) values and the newly inserted
# "05dec2025" (F) values for every data row, and fixes one unrelated data
# correction (row 18, column B: 6 -> 5).

# New E (04dec2025) / F (05dec2025) values, keyed by row number.
$newValues = @{
    2  = @(5, 0)
    3  = @(16, 0)
    4  = @(6, 0)
    5  = @(20, 0)
    6  = @(8, 0)
    7  = @(3, 0)
    8  = @(10, 0)
    9  = @(0, 8)
    10 = @(9, 0)
    11 = @(6, 0)
    12 = @(4, 0)
    13 = @(10, 0)
    14 = @(9, 11)
    15 = @(10, 0)
    16 = @(0, 0)
    17 = @(9, 0)
    18 = @(10, 0)
    19 = @(10, 0)
    20 = @(0, 7)
    21 = @(10, 0)
}

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("crosstab", "annot")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $isAnnot = ($sheetName -eq "annot")

    # Insert a new blank column at F; everything from F onward (26nov..30nov)
    # shifts right to G..K.
    $ws.Columns("F:F").Insert()

    # New header for the inserted column.
    if ($isAnnot) {
        $ws.Cells.Item(1, 6).Value = "05dec2025"
    } else {
        $ws.Cells.Item(1, 6).Value = "05dec2025"
    }

    foreach ($r in $newValues.Keys) {
        $pair = $newValues[$r]
        $eVal = $pair[0]
        $fVal = $pair[1]

        if ($isAnnot) {
            # annot sheet stores numbers as text, with 0 represented as a
            # truly empty cell.
            if ($eVal -eq 0) {
                $ws.Cells.Item($r, 5).Value = $null
            } else {
                $ws.Cells.Item($r, 5).Value = "$eVal"
            }
            if ($fVal -eq 0) {
                $ws.Cells.Item($r, 6).Value = $null
            } else {
                $ws.Cells.Item($r, 6).Value = "$fVal"
            }
        } else {
            $ws.Cells.Item($r, 5).Value = $eVal
            $ws.Cells.Item($r, 6).Value = $fVal
        }
    }

    # Unrelated data correction: row 18, column B, 6 -> 5.
    if ($isAnnot) {
        $ws.Cells.Item(18, 2).Value = "5"
    } else {
        $ws.Cells.Item(18, 2).Value = 5
    }
}
